{"js": "// The edit removes the placeholder run \"vnpt.SiteAddress\" that followed\n// the \"\u0110\u1ecba ch\u1ec9: \" label, leaving the label run (and the paragraph mark)\n// intact.\nconst body = context.document.body;\nconst results = body.search(\"vnpt.SiteAddress\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# The paragraph \"\u0110\u1ecba ch\u1ec9: \" was followed by a second run containing the\n# unresolved merge placeholder \"vnpt.SiteAddress\". This placeholder run\n# is removed entirely, leaving only the \"\u0110\u1ecba ch\u1ec9: \" label in that bullet.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"vnpt.SiteAddress\"\n$range.Find.Replacement.Text = \"\"\n$range.Find.Forward = $true\n$range.Find.Wrap = 1\n$range.Find.MatchCase = $true\n$range.Find.MatchWholeWord = $false\n\n$found = $range.Find.Execute()\nif ($range.Find.Found) {\n    $range.Delete()\n}\n"}
